$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: requirement id becomes text "CR1", and Comments gets a note
$ws.Range("A2").Value = "CR1"
$ws.Range("F2").Value = "This is demo"

# Remove the old "Simple" row (row 3); the former Total row (row 4) shifts up to row 3
$ws.Rows.Item(3).Delete()

# Update the (now) Total row totals to reflect removal of the "Simple" line
$ws.Range("D3").Value = 14.57
$ws.Range("E3").Value = 16.03
